$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.308.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.26%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.876.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.85%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.7123"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.43%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'242.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.76%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.01%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.78%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07761"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.09%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'24.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.10%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.08487"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +2.77%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.886.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.18%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'5.214"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.49%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.7104"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.08%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'91.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.24%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'29.311.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.31%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.000008257"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +5.85%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'6.007"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.37%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -0.79%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'2.132.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.22%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'13.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.65%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.9997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.00%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'7.822"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.97%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.00%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.1619"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.34%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'162.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.17%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.028"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.93%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'18.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.16%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.512"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.07%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.404"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.01%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'4.332"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +3.49%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -2.68%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.05238"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.85%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.931"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.99%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.178"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.51%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.7393"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.55%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.687"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.32%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +0.39%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.724"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.37%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'1.171.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.91%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'6.383"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +4.54%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'TrustWalletToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.8886"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.84%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'Aave"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'72.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.78%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'106.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +4.42%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +0.02%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.029.70"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.22%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.813"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.42%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.5206"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.46%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00000000122"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.16%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'9.399"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.68%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.4310"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.01%  "
$ws.Range("E51").Style = "Normal"
